# Apply Trade #9 closing update across the live trading results workbook.
$wb = $excel.ActiveWorkbook

$summary   = $wb.Worksheets.Item("Summary")
$strategy  = $wb.Worksheets.Item("Strategy Status")
$allTrades = $wb.Worksheets.Item("All Trades")
$mktMaking = $wb.Worksheets.Item("MarketMaking")

# --- Summary sheet: updated account-level stats ---
$summary.Range("B3").Value = 1200.05
$summary.Range("B4").Value = 0.05
$summary.Range("B5").Value = 0.11
$summary.Range("B6").Value = 9
$summary.Range("B7").Value = 4
$summary.Range("B9").Value = 44.44

# --- Strategy Status sheet: MarketMaking row (row 4) ---
$strategy.Range("C4").Value = 100.05
$strategy.Range("D4").Value = 9
$strategy.Range("E4").Value = 0.05
$strategy.Range("F4").Value = 0.05
$strategy.Range("G4").Value = 44.44

# --- All Trades & MarketMaking sheets: Trade #9 (row 10) closed out ---
foreach ($ws in @($allTrades, $mktMaking)) {
    $ws.Range("G10").Value = 0.858586
    $ws.Range("H10").Value = "CLOSED"
    $ws.Range("I10").Value = 3.4441
    $ws.Range("J10").Value = 0.03
    $ws.Range("K10").Value = 100.05
    $ws.Range("P10").Value = "early_exit"
    $ws.Range("Q10").Value = 0.13
}
